# Update the Cameroon MSME summary figures with more precise decimal values.
# The affected cells hold their numbers as text (shared strings), so we mark
# the range as Text before writing the new values to keep them stored as
# text rather than being auto-converted to numeric cells.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$range = $ws.Range("B10:D12")
$range.NumberFormat = "@"

# Enterprises (absolute #)
$ws.Range("B10").Value = "70476.75"
$ws.Range("C10").Value = "22552.56"
$ws.Range("D10").Value = "93029.31"

# Enterprises density (per 1000 people)
$ws.Range("B11").Value = "3.51"
$ws.Range("C11").Value = "1.12"
$ws.Range("D11").Value = "4.63"

# Employment (% of total)
$ws.Range("B12").Value = "23.46"
$ws.Range("C12").Value = "18.28"
$ws.Range("D12").Value = "41.74"
